$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 11 before reordering)
$ws.Range("A2").Value = 111406144
$ws.Range("AB2").Value = "16:48"
$ws.Range("AC2").Value = "Ett sextioplantor och tio blommor"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "60"
$ws.Range("Q2").Value = 564569.2468205031
$ws.Range("R2").Value = 6615734.803793807
$ws.Range("S2").Value = 4
$ws.Range("Z2").Value = "16:47"

# Row 3 (was row 2 before reordering)
$ws.Range("A3").Value = 111405420
$ws.Range("AB3").Value = "16:10"
$ws.Range("AC3").ClearContents()
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "30"
$ws.Range("Q3").Value = 564512.6192034025
$ws.Range("R3").Value = 6615790.408539454
$ws.Range("S3").Value = 7
$ws.Range("Z3").Value = "16:10"

# Row 4 (was row 12 before reordering)
$ws.Range("A4").Value = 111404914
$ws.Range("AB4").Value = "15:46"
$ws.Range("AC4").Value = "60 plantor och flera blommor."
$ws.Range("I4").Value = ""
$ws.Range("Q4").Value = 564468.1590201143
$ws.Range("R4").Value = 6615756.803563487
$ws.Range("Z4").Value = "15:46"

# Row 5 (was row 10 before reordering)
$ws.Range("A5").Value = 111405158
$ws.Range("AB5").Value = "16:00"
$ws.Range("AC5").ClearContents()
$ws.Range("B5").Value = 95524
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 221944
$ws.Range("F5").Value = "Lopplummer"
$ws.Range("G5").Value = "Huperzia selago"
$ws.Range("H5").Value = "(L.) Bernh. ex Schrank & Mart."
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("Q5").Value = 564459.9772311078
$ws.Range("R5").Value = 6615761.714472669
$ws.Range("S5").Value = 4
$ws.Range("Z5").Value = "16:00"

# Row 6 (was row 4 before reordering)
$ws.Range("A6").Value = 111407346
$ws.Range("AB6").Value = "17:51"
$ws.Range("AC6").Value = "Ca 80 plantor och 5 blommor"
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "80"
$ws.Range("Q6").Value = 564656.4028239732
$ws.Range("R6").Value = 6615816.131796785
$ws.Range("Z6").Value = "17:50"

# Row 7 (was row 9 before reordering)
$ws.Range("A7").Value = 111405082
$ws.Range("AB7").Value = "15:55"
$ws.Range("AC7").Value = "Ett femtiotal knärotsplantor. Tre blommor"
$ws.Range("Q7").Value = 564461.006485557
$ws.Range("R7").Value = 6615760.721820729
$ws.Range("S7").Value = 4
$ws.Range("Z7").Value = "15:54"

# Row 8 (was row 7 before reordering)
$ws.Range("A8").Value = 111405594
$ws.Range("AB8").Value = "16:16"
$ws.Range("AC8").Value = "Växer 5 m från planerad basväg."
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "50"
$ws.Range("Q8").Value = 564544.0294749426
$ws.Range("R8").Value = 6615788.424780905
$ws.Range("S8").Value = 10
$ws.Range("Z8").Value = "16:16"

# Row 9 (was row 5 before reordering)
$ws.Range("A9").Value = 111405823
$ws.Range("AB9").Value = "16:32"
$ws.Range("AC9").Value = "Mer än 150 plantor i blåbärsriset"
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "150"
$ws.Range("Q9").Value = 564551.2081450538
$ws.Range("R9").Value = 6615753.689971274
$ws.Range("S9").Value = 7
$ws.Range("Z9").Value = "16:32"

# Row 10 (was row 8 before reordering)
$ws.Range("A10").Value = 111405472
$ws.Range("AB10").Value = "16:13"
$ws.Range("AC10").Value = "En blomma."
$ws.Range("B10").Value = 96348
$ws.Range("D10").Value = "VU"
$ws.Range("E10").Value = 220787
$ws.Range("F10").Value = "Knärot"
$ws.Range("G10").Value = "Goodyera repens"
$ws.Range("H10").Value = "(L.) R. Br."
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = "7"
$ws.Range("J10").Value = "plantor/tuvor"
$ws.Range("K10").Value = "blomning"
$ws.Range("Q10").Value = 564524.4037030815
$ws.Range("R10").Value = 6615811.324952397
$ws.Range("S10").Value = 5
$ws.Range("Z10").Value = "16:13"

# Row 11 (was row 3 before reordering)
$ws.Range("A11").Value = 111407769
$ws.Range("AB11").Value = "18:08"
$ws.Range("AC11").Value = "Ett tjugotal plantor och tre blommor fem meter från planerad basväg."
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = "20"
$ws.Range("Q11").Value = 564530.9713830581
$ws.Range("R11").Value = 6615753.34056537
$ws.Range("Z11").Value = "18:07"

# Row 12 (was row 6 before reordering)
$ws.Range("A12").Value = 111405323
$ws.Range("AB12").Value = "16:06"
$ws.Range("AC12").Value = "Hundra plantor i mossan. Fyra blommor"
$ws.Range("Q12").Value = 564521.8477020637
$ws.Range("R12").Value = 6615783.495059335
$ws.Range("Z12").Value = "16:06"
